# "Delete excel datas remove" - clean up the Sayfa1 score sheet:
#  - remove the stray Yunanistan entry in C3 (cell content delete, not row delete)
#  - normalize the leftover numeric score cells
#  - move the active selection to C3

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

$ws.Range("C3").ClearContents()

$ws.Range("E2").Value = 85
$ws.Range("E5").Value = 78
$ws.Range("E10").Value = 72

$null = $ws.Range("C3").Select()
